$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "per ..." pallet-parameter sub-headers (row 2, columns J/K)
# "Per Package" -> "PCS Per PU" and "Per Handling Unit" -> "PU per HU"
$ws.Range("J2").Value = "PCS Per PU"
$ws.Range("K2").Value = "PU per HU"

# Move the current selection to K3 (matches the saved view state in the file)
$null = $ws.Range("K3").Select()
